$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Publication-years data, refreshed from a new search query.
# Column A holds the year as text (leading apostrophe keeps it text-typed,
# same as the original "number stored as text" cells); column B holds the count.
$years  = @("2025", "2024", "2023", "2022", "2021", "2020", "2019", "2018", "2017", "2016", "2015", "2014", "2013", "2012", "2011", "2010", "2009", "2008", "2007", "2006", "2005", "2004", "2003", "2002", "2001", "2000", "1999", "1998", "1997", "1996", "1995", "1994", "1993", "1992", "1991", "1989", "1988", "1987", "1986", "1983", "1979", "1969", "1964", "1954")
$counts = @(63, 2817, 2652, 2398, 2184, 1836, 1473, 1161, 981, 800, 636, 520, 437, 349, 306, 231, 172, 164, 146, 79, 66, 51, 49, 42, 38, 38, 24, 21, 18, 20, 19, 13, 4, 7, 5, 1, 1, 2, 1, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = "'" + $years[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

Write-Host "Rows written:" $years.Length
Write-Host "UsedRange:" $ws.UsedRange.Address()
